$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (veg_id renumbering + an extra CCl row + n_value tweaks) ---
$ws.Range("A2").Value = 100

$ws.Range("A3").Value = 101
$ws.Range("B3").Value = "CCl"
$ws.Range("C3").Value = "cottonwood closed"
$ws.Range("D3").Value = 0.15

$ws.Range("A4").Value = 102
$ws.Range("B4").Value = "WE"
$ws.Range("C4").Value = "willow established"

$ws.Range("A5").Value = 106
$ws.Range("B5").Value = "WC"
$ws.Range("C5").Value = "willow cottonwood"
$ws.Range("D5").Value = 0.15

# --- Insert two new rows (6, 7), carrying the formatting down from the row above ---
$ws.Rows.Item(6).Insert(-4121, 0)
$ws.Range("A6").Value = 210
$ws.Range("B6").Value = "CJ"
$ws.Range("C6").Value = "cottonwood juniper"
$ws.Range("D6").Value = 0.15

$ws.Rows.Item(7).Insert(-4121, 0)
$ws.Range("A7").Value = 215
$ws.Range("B7").Value = "CJ"
$ws.Range("C7").Value = "cottonwood juniper"
$ws.Range("D7").Value = 0.1

# --- Column widths for C and D ---
$ws.Columns.Item(3).ColumnWidth = 17.8307291666667
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection on D7 (last edited cell) ---
$ws.Range("D7").Select()
